# Updated cryptos list on Sun Sep 17 11:31:53 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row with
# newly-scraped values. A couple of numeric-looking price strings are
# written with a leading apostrophe so Excel keeps them as literal text
# (matching the original sheet, where prices are stored as text, not
# numbers) instead of silently parsing them into floating point numbers.
# Row 50/51 also swap which coin (USDD / Mantle) occupies which rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.843.27"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").Value = "1.646.43"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").Value = "'216.79"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").Value = "'0.0628"
$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("D10").Value = "'19.22"
$ws.Range("E10").Value = "  -0.18%  "

$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "1.643.05"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("E13").Value = "  -0.51%  "

$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").Value = "'64.76"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").Value = "26.850.31"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").Value = "'214.47"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").Value = "'4.39"
$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("E21").Value = "  +10.92%  "

$ws.Range("D22").Value = "'6.26"
$ws.Range("E22").Value = "  -0.51%  "

$ws.Range("D23").Value = "'9.35"
$ws.Range("E23").Value = "  -1.46%  "

$ws.Range("D24").Value = "'147.04"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("E25").Value = "  +0.41%  "

$ws.Range("D26").Value = "'0.119"
$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("D27").Value = "'7.20"
$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("D28").Value = "'15.70"
$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("E29").Value = "  -1.74%  "

$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("D31").Value = "'3.36"
$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").Value = "1.297.71"
$ws.Range("E33").Value = "  +1.54%  "

$ws.Range("D34").Value = "'1.53"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("E36").Value = "  -1.48%  "

$ws.Range("D37").Value = "'0.535"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("D38").Value = "'0.824"
$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("D40").Value = "'0.808"
$ws.Range("E40").Value = "  -1.00%  "

$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("E42").Value = "  -2.33%  "

$ws.Range("D43").Value = "1.788.23"
$ws.Range("E43").Value = "  +0.32%  "

$ws.Range("D44").Value = "'61.60"
$ws.Range("E44").Value = "  +2.92%  "

$ws.Range("D45").Value = "'92.01"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("D48").Value = "'7.69"
$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("D49").Value = "'0.0970"
$ws.Range("E49").Value = "  +0.04%  "

# Row 50 was USDD, now becomes Mantle
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.407"
$ws.Range("E50").Value = "  +0.12%  "

# Row 51 was Mantle, now becomes USDD
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.01"
$ws.Range("E51").Value = "  +0.51%  "
